# "Plying with the ammonia script, nothing major"
#
# The author appended 9 more experimental-data rows (164:172) below the
# existing table on Sheet1, reusing the same shared-formula pattern that
# was already being auto-filled down column D:I (P_SI, T_SI, rho_SI and
# the delta_* columns), and then left the selection sitting on G175.
#
# Columns: A=P_raw  B=T_raw  C=rho_raw  D=P_SI  E=T_SI  F=rho_SI
#          G=delta_P  H=delta_T  I=delta_rho

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New raw measurements (P_raw, T_raw, rho_raw) for rows 164-172, taken
# straight from the source data the rest of the sheet already follows.
$newRows = @(
    @(164, 160.82,  313.141, 676.92),
    @(165, 104.15,  313.141, 652.47),
    @(166,  28.685, 313.141, 605.88),
    @(167,  86.75,  317.939, 638.92),
    @(168,  56.05,  317.939, 620.46),
    @(169,  14.007, 317.939, 585.88),
    @(170,  50.662, 323.138, 610.28),
    @(171,  25.331, 323.138, 589.39),
    @(172,   5.066, 323.138, 567.2)
)

foreach ($r in $newRows) {
    $row = $r[0]

    # Raw inputs (columns A:C) - plain literals, same as every row above.
    $ws.Cells.Item($row, 1).Value2 = $r[1]
    $ws.Cells.Item($row, 2).Value2 = $r[2]
    $ws.Cells.Item($row, 3).Value2 = $r[3]

    # Derived SI columns (D:F) and the delta_* columns (G:I) - same
    # formulas as the shared-formula groups si="12".."17" used by every
    # row from 131 down to 163, just continued for the new rows.
    $ws.Cells.Item($row, 4).Formula = "=A$row*1000000"
    $ws.Cells.Item($row, 5).Formula = "=B$row"
    $ws.Cells.Item($row, 6).Formula = "=C$row"
    $ws.Cells.Item($row, 7).Formula = "=0.001*100"
    $ws.Cells.Item($row, 8).Formula = "=0.003/E$row*100"
    $ws.Cells.Item($row, 9).Formula = "=0.001*100"
}

# Leave the sheet scrolled/selected the way the author left it.
$ws.Range("G175").Select() | Out-Null

# Best-effort: also nudge the window scroll position to roughly where
# the new rows live (harmless if the host doesn't persist it).
try {
    $wb.Windows.Item(1).ScrollRow = 154
    $wb.Windows.Item(1).ScrollColumn = 1
} catch {
}
